$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 95.63567833333333
$ws.Range("H2").Value2 = 286.907035
$ws.Range("I2").Value2 = 0.2808828217467972
$ws.Range("J2").Value2 = 0.2808828217467972
$ws.Range("M2").Value2 = 0.05057900000000001
$ws.Range("N2").Value2 = 0.151737
$ws.Range("O2").Value2 = 0.01400296657613869
$ws.Range("P2").Value2 = 0.01400296657613869
$ws.Range("Q2").Value2 = 4.837156974421667
$ws.Range("R2").Value2 = 43.534412769795
$ws.Range("S2").Value2 = 0.003933192764731923
$ws.Range("T2").Value2 = 0.003933192764731924

# Row 3
$ws.Range("G3").Value2 = 95.63567833333333
$ws.Range("H3").Value2 = 286.907035
$ws.Range("I3").Value2 = 0.2808828217467972
$ws.Range("J3").Value2 = 0.2808828217467972
$ws.Range("O3").Value2 = 0.146324388539341
$ws.Range("P3").Value2 = 0.146324388539341
$ws.Range("Q3").Value2 = 50.54600628392222
$ws.Range("R3").Value2 = 454.9140565553
$ws.Range("S3").Value2 = 0.04110000714330481
$ws.Range("T3").Value2 = 0.04110000714330481

# Row 4
$ws.Range("G4").Value2 = 95.63567833333333
$ws.Range("H4").Value2 = 286.907035
$ws.Range("I4").Value2 = 0.2808828217467972
$ws.Range("J4").Value2 = 0.2808828217467972
$ws.Range("O4").Value2 = 0.8396726448845202
$ws.Range("P4").Value2 = 0.8396726448845202
$ws.Range("Q4").Value2 = 290.0548514737822
$ws.Range("R4").Value2 = 2610.49366326404
$ws.Range("S4").Value2 = 0.2358496218387604
$ws.Range("T4").Value2 = 0.2358496218387605

# Row 5
$ws.Range("I5").Value2 = 0.392628215788982
$ws.Range("J5").Value2 = 0.392628215788982
$ws.Range("M5").Value2 = 0.05057900000000001
$ws.Range("N5").Value2 = 0.151737
$ws.Range("O5").Value2 = 0.01400296657613869
$ws.Range("P5").Value2 = 0.01400296657613869
$ws.Range("Q5").Value2 = 6.761553805773334
$ws.Range("R5").Value2 = 60.85398425196001
$ws.Range("S5").Value2 = 0.005497959782542084
$ws.Range("T5").Value2 = 0.005497959782542084

# Row 6
$ws.Range("I6").Value2 = 0.392628215788982
$ws.Range("J6").Value2 = 0.392628215788982
$ws.Range("O6").Value2 = 0.146324388539341
$ws.Range("P6").Value2 = 0.146324388539341
$ws.Range("R6").Value2 = 635.8954002664
$ws.Range("S6").Value2 = 0.05745108359861522
$ws.Range("T6").Value2 = 0.05745108359861522

# Row 7
$ws.Range("I7").Value2 = 0.392628215788982
$ws.Range("J7").Value2 = 0.392628215788982
$ws.Range("O7").Value2 = 0.8396726448845202
$ws.Range("P7").Value2 = 0.8396726448845202
$ws.Range("S7").Value2 = 0.3296791724078246
$ws.Range("T7").Value2 = 0.3296791724078246

# Row 8
$ws.Range("I8").Value2 = 0.3264889624642208
$ws.Range("J8").Value2 = 0.3264889624642208
$ws.Range("M8").Value2 = 0.05057900000000001
$ws.Range("N8").Value2 = 0.151737
$ws.Range("O8").Value2 = 0.01400296657613869
$ws.Range("P8").Value2 = 0.01400296657613869
$ws.Range("Q8").Value2 = 5.622552322830001
$ws.Range("R8").Value2 = 50.60297090547
$ws.Range("S8").Value2 = 0.004571814028864683
$ws.Range("T8").Value2 = 0.004571814028864684

# Row 9
$ws.Range("I9").Value2 = 0.3264889624642208
$ws.Range("J9").Value2 = 0.3264889624642208
$ws.Range("O9").Value2 = 0.146324388539341
$ws.Range("P9").Value2 = 0.146324388539341
$ws.Range("S9").Value2 = 0.04777329779742095
$ws.Range("T9").Value2 = 0.04777329779742096

# Row 10
$ws.Range("I10").Value2 = 0.3264889624642208
$ws.Range("J10").Value2 = 0.3264889624642208
$ws.Range("O10").Value2 = 0.8396726448845202
$ws.Range("P10").Value2 = 0.8396726448845202
$ws.Range("S10").Value2 = 0.2741438506379351
$ws.Range("T10").Value2 = 0.2741438506379351
